# Update "想去人数" (interested-count) figures for a couple of events that
# appear on both the "展览" sheet and the consolidated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F4").Value = 182
    $ws.Range("F5").Value = 3361

    if ($name -eq "展览") {
        $ws.Range("F7").Value = 17
        $ws.Range("F8").Value = 422
    }
    elseif ($name -eq "全部类型") {
        $ws.Range("F9").Value = 17
        $ws.Range("F10").Value = 422
    }
}
